$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the existing header cell H1 into the new
# header cells I1 and J1, then set their text (I0 / IF).
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data for columns I and J (rows 2-47): row -> (I value, J value)
$newData = @(
    ,@(2, 9, 9)
    ,@(3, 8, 8)
    ,@(4, 8, 8)
    ,@(5, 7, 7)
    ,@(6, 9, 9)
    ,@(7, 6, 6)
    ,@(8, 9, 9)
    ,@(9, 6, 6)
    ,@(10, 6, 6)
    ,@(11, 7, 7)
    ,@(12, 8, 8)
    ,@(13, 7, 7)
    ,@(14, 7, 7)
    ,@(15, 6, 6)
    ,@(16, 7, 7)
    ,@(17, 7, 7)
    ,@(18, 8, 8)
    ,@(19, 8, 8)
    ,@(20, 8, 8)
    ,@(21, 9, 9)
    ,@(22, 6, 6)
    ,@(23, 9, 9)
    ,@(24, 7, 7)
    ,@(25, 8, 8)
    ,@(26, 8, 8)
    ,@(27, 9, 9)
    ,@(28, 7, 7)
    ,@(29, 8, 8)
    ,@(30, 7, 7)
    ,@(31, 8, 8)
    ,@(32, 6, 6)
    ,@(33, 9, 9)
    ,@(34, 9, 9)
    ,@(35, 7, 7)
    ,@(36, 8, 8)
    ,@(37, 7, 7)
    ,@(38, 8, 8)
    ,@(39, 9, 9)
    ,@(40, 7, 7)
    ,@(41, 9, 9)
    ,@(42, 7, 7)
    ,@(43, 8, 8)
    ,@(44, 5, 5)
    ,@(45, 6, 6)
    ,@(46, 5, 5)
    ,@(47, 4, 4)
)

foreach ($entry in $newData) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}

